$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.333947
$ws.Cells.Item(2, 8).Value = 1.001841
$ws.Cells.Item(2, 9).Value = 0.9184595666969813
$ws.Cells.Item(2, 10).Value = 0.9184595666969813
$ws.Cells.Item(2, 13).Value = 45.924193
$ws.Cells.Item(2, 14).Value = 137.772579
$ws.Cells.Item(2, 15).Value = 0.307792367338991
$ws.Cells.Item(2, 16).Value = 0.307792367338991
$ws.Cells.Item(2, 17).Value = 15.336246479771
$ws.Cells.Item(2, 18).Value = 138.026218317939
$ws.Cells.Item(2, 19).Value = 0.2826948443388078
$ws.Cells.Item(2, 20).Value = 0.2826948443388078

$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.333947
$ws.Cells.Item(3, 8).Value = 1.001841
$ws.Cells.Item(3, 9).Value = 0.9184595666969813
$ws.Cells.Item(3, 10).Value = 0.9184595666969813
$ws.Cells.Item(3, 15).Value = 0.2696759485354523
$ws.Cells.Item(3, 16).Value = 0.2696759485354523
$ws.Cells.Item(3, 17).Value = 13.43703501214733
$ws.Cells.Item(3, 18).Value = 120.933315109326
$ws.Cells.Item(3, 19).Value = 0.2476864548404689
$ws.Cells.Item(3, 20).Value = 0.2476864548404689

$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.333947
$ws.Cells.Item(4, 8).Value = 1.001841
$ws.Cells.Item(4, 9).Value = 0.9184595666969813
$ws.Cells.Item(4, 10).Value = 0.9184595666969813
$ws.Cells.Item(4, 13).Value = 23.60320766666667
$ws.Cells.Item(4, 14).Value = 70.809623
$ws.Cells.Item(4, 15).Value = 0.1581930283351338
$ws.Cells.Item(4, 16).Value = 0.1581930283351339
$ws.Cells.Item(4, 17).Value = 7.882220390660333
$ws.Cells.Item(4, 18).Value = 70.939983515943
$ws.Cells.Item(4, 19).Value = 0.1452939002591703
$ws.Cells.Item(4, 20).Value = 0.1452939002591703

$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.333947
$ws.Cells.Item(5, 8).Value = 1.001841
$ws.Cells.Item(5, 9).Value = 0.9184595666969813
$ws.Cells.Item(5, 10).Value = 0.9184595666969813
$ws.Cells.Item(5, 13).Value = 19.226538
$ws.Cells.Item(5, 14).Value = 57.679614
$ws.Cells.Item(5, 15).Value = 0.1288597852280838
$ws.Cells.Item(5, 16).Value = 0.1288597852280838
$ws.Cells.Item(5, 17).Value = 6.420644685486001
$ws.Cells.Item(5, 18).Value = 57.785802169374
$ws.Cells.Item(5, 19).Value = 0.1183525025052519
$ws.Cells.Item(5, 20).Value = 0.1183525025052519

$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.333947
$ws.Cells.Item(6, 8).Value = 1.001841
$ws.Cells.Item(6, 9).Value = 0.9184595666969813
$ws.Cells.Item(6, 10).Value = 0.9184595666969813
$ws.Cells.Item(6, 13).Value = 20.21413933333333
$ws.Cells.Item(6, 14).Value = 60.64241799999999
$ws.Cells.Item(6, 15).Value = 0.1354788705623391
$ws.Cells.Item(6, 16).Value = 0.1354788705623391
$ws.Cells.Item(6, 17).Value = 6.750451187948666
$ws.Cells.Item(6, 18).Value = 60.75406069153799
$ws.Cells.Item(6, 19).Value = 0.1244318647532824
$ws.Cells.Item(6, 20).Value = 0.1244318647532824

$ws.Cells.Item(7, 9).Value = 0.08154043330301874
$ws.Cells.Item(7, 10).Value = 0.08154043330301874
$ws.Cells.Item(7, 13).Value = 45.924193
$ws.Cells.Item(7, 14).Value = 137.772579
$ws.Cells.Item(7, 15).Value = 0.307792367338991
$ws.Cells.Item(7, 16).Value = 0.307792367338991
$ws.Cells.Item(7, 17).Value = 1.361545165999667
$ws.Cells.Item(7, 18).Value = 12.253906493997
$ws.Cells.Item(7, 19).Value = 0.02509752300018324
$ws.Cells.Item(7, 20).Value = 0.02509752300018324

$ws.Cells.Item(8, 9).Value = 0.08154043330301874
$ws.Cells.Item(8, 10).Value = 0.08154043330301874
$ws.Cells.Item(8, 15).Value = 0.2696759485354523
$ws.Cells.Item(8, 16).Value = 0.2696759485354523
$ws.Cells.Item(8, 19).Value = 0.02198949369498336
$ws.Cells.Item(8, 20).Value = 0.02198949369498336

$ws.Cells.Item(9, 9).Value = 0.08154043330301874
$ws.Cells.Item(9, 10).Value = 0.08154043330301874
$ws.Cells.Item(9, 13).Value = 23.60320766666667
$ws.Cells.Item(9, 14).Value = 70.809623
$ws.Cells.Item(9, 15).Value = 0.1581930283351338
$ws.Cells.Item(9, 16).Value = 0.1581930283351339
$ws.Cells.Item(9, 17).Value = 0.6997800331654445
$ws.Cells.Item(9, 18).Value = 6.298020298489
$ws.Cells.Item(9, 19).Value = 0.01289912807596353
$ws.Cells.Item(9, 20).Value = 0.01289912807596354

$ws.Cells.Item(10, 9).Value = 0.08154043330301874
$ws.Cells.Item(10, 10).Value = 0.08154043330301874
$ws.Cells.Item(10, 13).Value = 19.226538
$ws.Cells.Item(10, 14).Value = 57.679614
$ws.Cells.Item(10, 15).Value = 0.1288597852280838
$ws.Cells.Item(10, 16).Value = 0.1288597852280838
$ws.Cells.Item(10, 17).Value = 0.570021989778
$ws.Cells.Item(10, 18).Value = 5.130197908002
$ws.Cells.Item(10, 19).Value = 0.01050728272283189
$ws.Cells.Item(10, 20).Value = 0.01050728272283189

$ws.Cells.Item(11, 9).Value = 0.08154043330301874
$ws.Cells.Item(11, 10).Value = 0.08154043330301874
$ws.Cells.Item(11, 13).Value = 20.21413933333333
$ws.Cells.Item(11, 14).Value = 60.64241799999999
$ws.Cells.Item(11, 15).Value = 0.1354788705623391
$ws.Cells.Item(11, 16).Value = 0.1354788705623391
$ws.Cells.Item(11, 17).Value = 0.5993020649082221
$ws.Cells.Item(11, 18).Value = 5.393718584173999
$ws.Cells.Item(11, 19).Value = 0.01104700580905672
$ws.Cells.Item(11, 20).Value = 0.01104700580905672
